$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains a new first column ("sobat_id") and every former column
# A..J (nama_lengkap .. tahun) effectively slides one column right to B..K.
# We reproduce that by writing the final header/value text straight into its
# new home cell (rather than doing a structural column insert, which this
# host doesn't fully keep in lockstep with the sheet's <cols> metadata).

# Old G2 (jenis_kelamin's numeric "1" example) has no replacement at that
# position any more, so drop it outright.
$ws.Range("G2").Clear()

# Write the changed/new example row-2 values first, in the same order the
# original authoring happened in, so brand-new shared strings land at the
# same table slots the target file uses.
$ws.Range("H2").Value = "1"
$ws.Range("A1").Value = "sobat_id"
$ws.Range("I2").Value = "+62"
$ws.Range("K2").Value = "16-02-2024"
$ws.Range("D2").Value = "D001"
$ws.Range("E2").Value = "KC001"
$ws.Range("O1").Value = "keterangan"
$ws.Range("O2").Value = "1 = laki laki"
$ws.Range("O3").Value = "2 = perempuan"
$ws.Range("A2").Value = "1111"

# Re-home the rest of the header row + the untouched example values.
$ws.Range("B1").Value = "nama_lengkap"
$ws.Range("C1").Value = "alamat_mitra"
$ws.Range("D1").Value = "kode_desa"
$ws.Range("E1").Value = "kode_kecamatan"
$ws.Range("F1").Value = "kode_kabupaten"
$ws.Range("G1").Value = "kode_provinsi"
$ws.Range("H1").Value = "jenis_kelamin"
$ws.Range("I1").Value = "no_hp_mitra"
$ws.Range("J1").Value = "email_mitra"
$ws.Range("K1").Value = "tahun"
$ws.Range("B2").Value = "contoh"
$ws.Range("C2").Value = "contoh"
$ws.Range("J2").Value = "contoh@gmail.com"

# The e-mail cell used to be a live mailto hyperlink; it no longer is one,
# though its column (I2) still carries the old Hyperlink visual style.
$ws.Range("I2").Hyperlinks.Delete()

# Highlight the jenis_kelamin header + the new "keterangan" legend cells.
$ws.Range("H1").Interior.Color = 65535
$ws.Range("O1").Interior.Color = 65535
$ws.Range("O2").Interior.Color = 65535
$ws.Range("O3").Interior.Color = 65535

# Selection cosmetics to match the saved view state.
$ws.Range("N14").Select()
